# Poland Ekstraklasa workbook update - 27-03-2024 20:23
# Updates the two most-recently-finished fixtures (rows 223-224, ids 221-222) with final
# match data (previously placeholder rows with no score), inserts the fixture that was
# previously the last placeholder row as a finished match (new row 225, id 223), and
# appends 9 new upcoming fixtures (rows 226-234, ids 224-232) with opening odds only.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colNum = @{'A'=1; 'B'=2; 'C'=3; 'D'=4; 'E'=5; 'F'=6; 'G'=7; 'H'=8; 'I'=9; 'J'=10; 'K'=11; 'L'=12; 'M'=13; 'N'=14; 'O'=15; 'P'=16; 'Q'=17; 'R'=18; 'S'=19; 'T'=20; 'U'=21; 'V'=22; 'W'=23; 'X'=24; 'Y'=25; 'Z'=26; 'AA'=27; 'AB'=28; 'AC'=29}

# New rows need A (id) and E (Date) cell formatting (bordered/centered int style, date
# number format) copied from an existing, untouched data row (row 222) so the workbook
# keeps reusing the same two style records instead of minting new ones.
$styleSrcRow = 222
$newRows = @(225,226,227,228,229,230,231,232,233,234)
foreach ($r in $newRows) {
    $ws.Cells.Item($styleSrcRow, $colNum['A']).Copy() | Out-Null
    $ws.Cells.Item($r, $colNum['A']).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($styleSrcRow, $colNum['E']).Copy() | Out-Null
    $ws.Cells.Item($r, $colNum['E']).PasteSpecial(-4122) | Out-Null
}

# --- Row data (column letter => value) ---
$rowData = @{}
$rowData[223] = [ordered]@{
    'A' = 221
    'B' = 6775569
    'C' = "Poland Ekstraklasa"
    'D' = "Poland Ekstraklasa"
    'E' = 45368.35416666666
    'F' = "Korona Kielce"
    'G' = "Pogon Szczecin"
    'H' = 2
    'I' = 2
    'J' = "D"
    'K' = 3.4
    'L' = 3.2
    'M' = 2.2
    'N' = 3.3
    'O' = 3.25
    'P' = 2.25
    'Q' = 0.25
    'R' = 1.95
    'S' = 1.9
    'T' = 2.5
    'U' = 1.975
    'V' = 1.875
    'W' = -1
    'X' = 2.25
    'Y' = -1
    'Z' = 0.475
    'AA' = -0.5
    'AB' = 0.9750000000000001
    'AC' = -1
}
$rowData[224] = [ordered]@{
    'A' = 222
    'B' = 6774467
    'C' = "Poland Ekstraklasa"
    'D' = "Poland Ekstraklasa"
    'E' = 45368.45833333334
    'F' = "LKS Lodz"
    'G' = "Rakow Czestochowa"
    'H' = 1
    'I' = 1
    'J' = "D"
    'K' = 5.5
    'L' = 4.333
    'M' = 1.533
    'N' = 6.5
    'O' = 4.75
    'P' = 1.444
    'Q' = 1.25
    'R' = 1.85
    'S' = 2
    'T' = 2.5
    'U' = 1.825
    'V' = 2.025
    'W' = -1
    'X' = 3.75
    'Y' = -1
    'Z' = 0.8500000000000001
    'AA' = -1
    'AB' = -1
    'AC' = 1.025
}
$rowData[225] = [ordered]@{
    'A' = 223
    'B' = 6775571
    'C' = "Poland Ekstraklasa"
    'D' = "Poland Ekstraklasa"
    'E' = 45368.5625
    'F' = "Legia Warsaw"
    'G' = "Piast Gliwice"
    'H' = 3
    'I' = 1
    'J' = "H"
    'K' = 1.75
    'L' = 3.5
    'M' = 4.75
    'N' = 1.571
    'O' = 3.75
    'P' = 6
    'Q' = -0.75
    'R' = 1.75
    'S' = 2.05
    'T' = 2.25
    'U' = 2.05
    'V' = 1.8
    'W' = 0.571
    'X' = -1
    'Y' = -1
    'Z' = 0.75
    'AA' = -1
    'AB' = 1.05
    'AC' = -1
}
$rowData[226] = [ordered]@{
    'A' = 224
    'B' = 6774469
    'C' = "Poland Ekstraklasa"
    'D' = "Poland Ekstraklasa"
    'E' = 45381.35416666666
    'F' = "Rakow Czestochowa"
    'G' = "Ruch Chorzow"
    'K' = 1.4
    'L' = 4.75
    'M' = 7.5
    'N' = 1.4
    'O' = 4.75
    'P' = 7.5
    'Q' = -1.25
    'R' = 1.925
    'S' = 1.925
    'T' = 2.75
    'U' = 2.025
    'V' = 1.825
    'W' = 0
    'X' = 0
    'Y' = 0
    'Z' = 0
    'AA' = 0
}
$rowData[227] = [ordered]@{
    'A' = 225
    'B' = 6774468
    'C' = "Poland Ekstraklasa"
    'D' = "Poland Ekstraklasa"
    'E' = 45381.45833333334
    'F' = "Jagiellonia Bialystok"
    'G' = "LKS Lodz"
    'K' = 1.4
    'L' = 5
    'M' = 7
    'N' = 1.4
    'O' = 5
    'P' = 7
    'Q' = -1.25
    'R' = 1.875
    'S' = 1.975
    'T' = 2.75
    'U' = 1.8
    'V' = 2.05
    'W' = 0
    'X' = 0
    'Y' = 0
    'Z' = 0
    'AA' = 0
}
$rowData[228] = [ordered]@{
    'A' = 226
    'B' = 6775574
    'C' = "Poland Ekstraklasa"
    'D' = "Poland Ekstraklasa"
    'E' = 45381.5625
    'F' = "Piast Gliwice"
    'G' = "Slask Wroclaw"
    'K' = 2.1
    'L' = 3.1
    'M' = 4
    'N' = 2.1
    'O' = 3.1
    'P' = 4
    'Q' = -0.25
    'R' = 1.8
    'S' = 2.05
    'T' = 2
    'U' = 1.975
    'V' = 1.875
    'W' = 0
    'X' = 0
    'Y' = 0
    'Z' = 0
    'AA' = 0
}
$rowData[229] = [ordered]@{
    'A' = 227
    'B' = 6775575
    'C' = "Poland Ekstraklasa"
    'D' = "Poland Ekstraklasa"
    'E' = 45381.66666666666
    'F' = "Pogon Szczecin"
    'G' = "Cracovia Krakow"
    'K' = 1.909
    'L' = 3.6
    'M' = 4
    'N' = 1.909
    'O' = 3.6
    'P' = 4
    'Q' = -0.5
    'R' = 1.925
    'S' = 1.925
    'T' = 2.5
    'U' = 1.8
    'V' = 2.05
    'W' = 0
    'X' = 0
    'Y' = 0
    'Z' = 0
    'AA' = 0
}
$rowData[230] = [ordered]@{
    'A' = 228
    'B' = 6774877
    'C' = "Poland Ekstraklasa"
    'D' = "Poland Ekstraklasa"
    'E' = 45383.3125
    'F' = "Puszcza Niepolomice"
    'G' = "Radomiak Radom"
    'K' = 2.625
    'L' = 3.4
    'M' = 2.6
    'N' = 2.625
    'O' = 3.4
    'P' = 2.6
    'Q' = 0
    'R' = 1.975
    'S' = 1.875
    'T' = 2.25
    'U' = 1.85
    'V' = 2
    'W' = 0
    'X' = 0
    'Y' = 0
    'Z' = 0
    'AA' = 0
}
$rowData[231] = [ordered]@{
    'A' = 229
    'B' = 6775576
    'C' = "Poland Ekstraklasa"
    'D' = "Poland Ekstraklasa"
    'E' = 45383.41666666666
    'F' = "Stal Mielec"
    'G' = "Lech Poznan"
    'K' = 4.333
    'L' = 3.4
    'M' = 1.85
    'N' = 4.333
    'O' = 3.4
    'P' = 1.85
    'Q' = 0.5
    'R' = 2
    'S' = 1.85
    'T' = 2.25
    'U' = 1.925
    'V' = 1.925
    'W' = 0
    'X' = 0
    'Y' = 0
    'Z' = 0
    'AA' = 0
}
$rowData[232] = [ordered]@{
    'A' = 230
    'B' = 6775578
    'C' = "Poland Ekstraklasa"
    'D' = "Poland Ekstraklasa"
    'E' = 45383.52083333334
    'F' = "Widzew Lodz"
    'G' = "Korona Kielce"
    'K' = 2.25
    'L' = 3.2
    'M' = 3.4
    'N' = 2.25
    'O' = 3.2
    'P' = 3.4
    'Q' = -0.25
    'R' = 1.925
    'S' = 1.925
    'T' = 2.5
    'U' = 2.025
    'V' = 1.825
    'W' = 0
    'X' = 0
    'Y' = 0
    'Z' = 0
    'AA' = 0
}
$rowData[233] = [ordered]@{
    'A' = 231
    'B' = 6775573
    'C' = "Poland Ekstraklasa"
    'D' = "Poland Ekstraklasa"
    'E' = 45383.625
    'F' = "Gornik Zabrze"
    'G' = "Legia Warsaw"
    'K' = 3.6
    'L' = 3.5
    'M' = 2.05
    'N' = 3.8
    'O' = 3.5
    'P' = 1.95
    'Q' = 0.5
    'R' = 1.875
    'S' = 1.975
    'T' = 2.5
    'U' = 1.9
    'V' = 1.95
    'W' = 0
    'X' = 0
    'Y' = 0
    'Z' = 0
    'AA' = 0
}
$rowData[234] = [ordered]@{
    'A' = 232
    'B' = 6775577
    'C' = "Poland Ekstraklasa"
    'D' = "Poland Ekstraklasa"
    'E' = 45384.58333333334
    'F' = "Warta Poznan"
    'G' = "Zaglebie Lubin"
    'K' = 2.9
    'L' = 3.1
    'M' = 2.55
    'N' = 2.9
    'O' = 3.1
    'P' = 2.55
    'Q' = 0
    'R' = 2.05
    'S' = 1.8
    'T' = 2.25
    'U' = 2.05
    'V' = 1.8
    'W' = 0
    'X' = 0
    'Y' = 0
    'Z' = 0
    'AA' = 0
}

foreach ($r in $rowData.Keys) {
    foreach ($col in $rowData[$r].Keys) {
        $ws.Cells.Item([int]$r, $colNum[$col]).Value = $rowData[$r][$col]
    }
}
